$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph "Stock Prices - Crime Rates - Airport Traffic ..." (Data Sets)
# Target run layout:
#   "Stock Prices "
#   "(Yahoo Finance"
#   [bookmark _GoBack]
#   ") "
#   "- Crime Rates - Airport Traffic "
#   "- areas affected by government employment"
#   " - approval ratings for the President, the house & senate"
# ------------------------------------------------------------------

# Find the insertion point right after "Stock Prices "
$rSP = $d.Content
$rSP.Find.Execute("Stock Prices ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterSP = $rSP.End

# Insert the new text "(Yahoo Finance) " right after "Stock Prices ". This
# touches/merges the whole paragraph's runs, so all the wall bookmarks
# protecting the untouched trailing text are added AFTER this point, using
# positions computed relative to the (now stable, since nothing to their
# right moved) original boundaries.
$ins = $d.Range($posAfterSP, $posAfterSP)
$ins.InsertAfter("(Yahoo Finance) ")

# Re-split "Stock Prices " from "(Yahoo Finance) ..." with a temporary wall
# placed exactly at the point we just inserted at.
$d.Bookmarks.Add("ZZ_WallSP", $d.Range($posAfterSP, $posAfterSP))

# Protect the two untouched trailing runs ("- areas affected..." and
# " - approval ratings...") with temporary wall bookmarks placed exactly at
# their existing (unmoved) run boundaries.
$rAreas = $d.Content
$rAreas.Find.Execute("– areas affected by government employment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posBeforeAreas = $rAreas.Start
$d.Bookmarks.Add("ZZ_WallAreas", $d.Range($posBeforeAreas, $posBeforeAreas))

$rApproval = $d.Content
$rApproval.Find.Execute(" – approval ratings for the President", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posBeforeApproval = $rApproval.Start
$d.Bookmarks.Add("ZZ_WallApproval", $d.Range($posBeforeApproval, $posBeforeApproval))

# The real bookmark _GoBack belongs between "(Yahoo Finance" and ") "
$bmPos = $posAfterSP + 14   # length of "(Yahoo Finance"

# Remove the old _GoBack bookmark (currently sitting inside "Kag|gle" below)
$d.Bookmarks("_GoBack").Delete()

# Re-create _GoBack at its new location
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# Split ") " from "- Crime Rates - Airport Traffic " with another temp wall
$wallPos2 = $bmPos + 2   # length of ") "
$d.Bookmarks.Add("ZZ_WallCrime", $d.Range($wallPos2, $wallPos2))

# Clean up the temporary wall bookmarks - the run splits they created remain
# in place even after the bookmark markers themselves are removed.
$d.Bookmarks("ZZ_WallAreas").Delete()
$d.Bookmarks("ZZ_WallApproval").Delete()
$d.Bookmarks("ZZ_WallSP").Delete()
$d.Bookmarks("ZZ_WallCrime").Delete()

# ------------------------------------------------------------------
# Paragraph "We plan to pull these data sets from Kaggle, Data.gov, ..."
# Target run layout:
#   "We plan to pull these data sets from Kaggle, Data.gov, and the "
#   "N"
#   "asdaq's historical data"
# ------------------------------------------------------------------

# Protect "N" / "asdaq's historical data" from being absorbed by the merge
# below with a temporary wall at their existing boundary.
$rNas = $d.Content
$rNas.Find.Execute("Nasdaq", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posBeforeN = $rNas.Start
$d.Bookmarks.Add("ZZ_WallNasdaq", $d.Range($posBeforeN, $posBeforeN))

# Merge "Kag" and "gle, Data.gov, and the " into a single run by inserting
# (and immediately removing) a throwaway character right at the boundary -
# this forces the engine to coalesce the two adjacent, identically
# formatted runs.
$rKag = $d.Content
$rKag.Find.Execute("Kag", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$posAfterKag = $rKag.End
$d.Range($posAfterKag, $posAfterKag).InsertAfter("Z")

$rKagZ = $d.Content
$rKagZ.Find.Execute("KagZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$zEnd = $rKagZ.End
$d.Range($zEnd - 1, $zEnd).Delete()

# Remove the temporary wall now that the merge is done.
$d.Bookmarks("ZZ_WallNasdaq").Delete()

Write-Output "Data Sets paragraph: $($d.Paragraphs(15).Range.Text)"
Write-Output "Pull paragraph: $($d.Paragraphs(16).Range.Text)"
